# member payment verification page object class
# Adds a new "MemberPayment" worksheet with a transaction verification
# data row, and updates the selection/active-tab state left over from
# the editing session.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # MemberDetails
$ws2 = $wb.Worksheets.Item(2)   # UserAccountTypes

# --- Update leftover selections on the existing sheets -------------------
$ws1.Range("B6").Select()
$ws2.Range("A1:A1048576").Select()

# --- Add the new MemberPayment sheet at the end of the workbook ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "MemberPayment"

# Write cell values in the same order the original author entered them so
# the shared-string table is rebuilt with a matching append order.
$ws3.Range("A2").Value = "TXN001"
$ws3.Range("C1").Value = "Transaction_Name"
$ws3.Range("B1").Value = "Transaction_Login"
$ws3.Range("E1").Value = "Transaction_Type"
$ws3.Range("D1").Value = "Transaction_Amount"
$ws3.Range("F1").Value = "Transaction_Description"
$ws3.Range("E2").Value = "Savings to Current"
$ws3.Range("F2").Value = "Miscellaneous Expense"
$ws3.Range("B2").Value = "TestUser01"
$ws3.Range("A1").Value = "DATA_SET_ID"
$ws3.Range("C2").Value = "TestUser01"
$ws3.Range("D2").Value = 534

# Number formatting for the transaction amount cell
$ws3.Range("D2").NumberFormat = "0.00"

# Header row styling - copy the header format already used on the other
# sheets (bold white font, blue fill, medium border, centered)
$ws1.Range("A1").Copy()
$ws3.Range("A1:F1").PasteSpecial(-4122)

# Column widths (matches widths used on the comparable sheet1 columns)
$ws3.Columns.Item(1).ColumnWidth = 12.7109375
$ws3.Columns.Item(2).ColumnWidth = 18.140625
$ws3.Columns.Item(3).ColumnWidth = 19
$ws3.Columns.Item(4).ColumnWidth = 18.85546875
$ws3.Columns.Item(5).ColumnWidth = 16.5703125
$ws3.Columns.Item(6).ColumnWidth = 22.5703125

# Row height for the header row
$ws3.Rows.Item(1).RowHeight = 15.75

# Select B1 and make this the active sheet/tab
$ws3.Range("B1").Select()
$ws3.Activate()

Write-Output "MemberPayment sheet added"
